# Updates cryptos list values (Price + Volume(1h) columns) plus the
# Aave / FraxShare row swap (ranks 46/47), per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "43.305.38", "112.67") that must
# stay plain text, matching the source data (t="inlineStr"). Force the
# NumberFormat to Text ("@") before assigning so COM does not coerce the
# string into a real number.
$dCells = @(
    'D2',
    'D3',
    'D5',
    'D6',
    'D7',
    'D9',
    'D10',
    'D11',
    'D12',
    'D14',
    'D15',
    'D16',
    'D17',
    'D18',
    'D20',
    'D21',
    'D22',
    'D23',
    'D24',
    'D25',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D37',
    'D38',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D48',
    'D49',
    'D50'
)
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.305.38'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '2.275.87'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('D5').Value = '112.67'
$ws.Range('E5').Value = '  -4.48%  '
$ws.Range('D6').Value = '264.81'
$ws.Range('E6').Value = '  -2.22%  '
$ws.Range('D7').Value = '0.621'
$ws.Range('E7').Value = '  -1.06%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.608'
$ws.Range('E9').Value = '  -3.14%  '
$ws.Range('D10').Value = '47.75'
$ws.Range('E10').Value = '  -4.70%  '
$ws.Range('D11').Value = '0.0931'
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').Value = '8.82'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('D14').Value = '15.51'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').Value = '2.619.38'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').Value = '0.854'
$ws.Range('E16').Value = '  -1.64%  '
$ws.Range('D17').Value = '2.274.79'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').Value = '43.221.53'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('E19').Value = '  -2.80%  '
$ws.Range('D20').Value = '6.79'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('D21').Value = '71.41'
$ws.Range('E21').Value = '  -1.99%  '
$ws.Range('D22').Value = '2.51'
$ws.Range('E22').Value = '  -1.01%  '
$ws.Range('D23').Value = '232.09'
$ws.Range('E23').Value = '  -1.18%  '
$ws.Range('D24').Value = '9.65'
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('D25').Value = '2.87'
$ws.Range('E25').Value = '  -2.14%  '
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('D27').Value = '11.32'
$ws.Range('E27').Value = '  -1.97%  '
$ws.Range('D28').Value = '3.91'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').Value = '40.40'
$ws.Range('E29').Value = '  -7.71%  '
$ws.Range('D30').Value = '3.34'
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('D31').Value = '2.24'
$ws.Range('E31').Value = '  -1.22%  '
$ws.Range('D32').Value = '172.01'
$ws.Range('E32').Value = '  -3.57%  '
$ws.Range('D33').Value = '21.32'
$ws.Range('E33').Value = '  -3.11%  '
$ws.Range('D34').Value = '0.0908'
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('D35').Value = '5.74'
$ws.Range('E35').Value = '  +2.31%  '
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('D37').Value = '4.65'
$ws.Range('E37').Value = '  -2.80%  '
$ws.Range('D38').Value = '0.0353'
$ws.Range('E38').Value = '  -1.71%  '
$ws.Range('E39').Value = '  -4.69%  '
$ws.Range('E40').Value = '  -6.63%  '
$ws.Range('D41').Value = '2.63'
$ws.Range('E41').Value = '  +9.58%  '
$ws.Range('D42').Value = '76.47'
$ws.Range('E42').Value = '  +5.10%  '
$ws.Range('D43').Value = '13.85'
$ws.Range('E43').Value = '  +4.44%  '
$ws.Range('D44').Value = '0.237'
$ws.Range('E44').Value = '  -4.13%  '
$ws.Range('D45').Value = '6.13'
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '104.27'
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '8.62'
$ws.Range('E49').Value = '  -2.64%  '
$ws.Range('D50').Value = '0.0991'
$ws.Range('E50').Value = '  -1.77%  '
$ws.Range('E51').Value = '  +0.73%  '
